$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header style (from H1, the last header cell) onto the
# two new header cells so they match the other headers exactly (bold,
# centered, bordered).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
